$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 15 (Caso -173, "PACHECO DE MELO J A /ALT/ 2300") entirely;
# all rows below shift up by one.
$ws.Rows.Item(15).Delete()
